$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B93").Value = "Banco Pichincha"
$ws.Range("C93").Value = "Avenida Amazonas N35-211 y Japon, Quito, Ecuador"
$ws.Range("D93").Value = "Banking"

$ws.Range("C96").Value = "Av. Libertadores, Monterrey, Mexico"

$ws.Range("B92").Value = "PetroEcuador"
$ws.Range("C92").Value = "Alpallana E8-86 Y Av. , 6 de Diciembre, Quito, Ecuador"

$ws.Range("C92").Select()
$excel.ActiveWindow.ScrollRow = 74
$excel.ActiveWindow.ScrollColumn = 1
